$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.302.79"
$ws.Range("E2").Value = "  -5.02%  "

$ws.Range("D3").Value = "1.560.45"
$ws.Range("E3").Value = "  -5.22%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").Value = "289.68"
$ws.Range("E6").Value = "  -3.67%  "

$ws.Range("D7").Value = "0.3720"
$ws.Range("E7").Value = "  -1.90%  "

$ws.Range("D8").Value = "49.27"
$ws.Range("E8").Value = "  -2.55%  "

$ws.Range("D9").Value = "0.3399"
$ws.Range("E9").Value = "  -2.87%  "

$ws.Range("D10").Value = "1.162"
$ws.Range("E10").Value = "  -4.57%  "

$ws.Range("D11").Value = "0.07628"
$ws.Range("E11").Value = "  -5.48%  "

$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").Value = "21.37"
$ws.Range("E13").Value = "  -3.39%  "

$ws.Range("D14").Value = "6.028"
$ws.Range("E14").Value = "  -4.40%  "

$ws.Range("D15").Value = "6.913"
$ws.Range("E15").Value = "  -4.69%  "

$ws.Range("D16").Value = "1.562.69"

$ws.Range("D17").Value = "0.00001126"
$ws.Range("E17").Value = "  -7.08%  "

$ws.Range("D18").Value = "89.78"
$ws.Range("E18").Value = "  -5.68%  "

$ws.Range("D19").Value = "0.06724"
$ws.Range("E19").Value = "  -3.73%  "

$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").Value = "6.224"
$ws.Range("E21").Value = "  -6.03%  "

$ws.Range("D22").Value = "16.53"
$ws.Range("E22").Value = "  -5.13%  "

$ws.Range("D23").Value = "0.5283"
$ws.Range("E23").Value = "  -7.61%  "

$ws.Range("D24").Value = "11.95"
$ws.Range("E24").Value = "  -3.88%  "

$ws.Range("D25").Value = "22.306.25"
$ws.Range("E25").Value = "  -5.00%  "

$ws.Range("D26").Value = "2.405"
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").Value = "2.814"
$ws.Range("E27").Value = "  -5.41%  "

$ws.Range("D28").Value = "20.16"
$ws.Range("E28").Value = "  -3.99%  "

$ws.Range("D29").Value = "146.05"
$ws.Range("E29").Value = "  -3.64%  "

$ws.Range("D30").Value = "4.973"
$ws.Range("E30").Value = "  -4.29%  "

$ws.Range("D31").Value = "125.29"
$ws.Range("E31").Value = "  -4.75%  "

$ws.Range("D32").Value = "1.733.59"
$ws.Range("E32").Value = "  -5.34%  "

$ws.Range("D33").Value = "6.167"
$ws.Range("E33").Value = "  -9.73%  "

$ws.Range("D34").Value = "1.006"
$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("D35").Value = "2.007"
$ws.Range("E35").Value = "  -6.23%  "

$ws.Range("D36").Value = "10.01"
$ws.Range("E36").Value = "  -10.36%  "

$ws.Range("D37").Value = "0.08486"
$ws.Range("E37").Value = "  -3.47%  "

$ws.Range("D38").Value = "0.02537"
$ws.Range("E38").Value = "  -5.66%  "

$ws.Range("D39").Value = "0.2306"
$ws.Range("E39").Value = "  -4.61%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "5.505"
$ws.Range("E40").Value = "  -6.90%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.319"
$ws.Range("E41").Value = "  +1.77%  "

$ws.Range("D42").Value = "0.06389"
$ws.Range("E42").Value = "  -5.86%  "

$ws.Range("D43").Value = "11.67"
$ws.Range("E43").Value = "  -9.17%  "

$ws.Range("D44").Value = "0.6330"
$ws.Range("E44").Value = "  -7.93%  "

$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "14.05"
$ws.Range("E46").Value = "  -9.57%  "

$ws.Range("D47").Value = "0.5960"
$ws.Range("E47").Value = "  -6.57%  "

$ws.Range("D48").Value = "3.755"
$ws.Range("E48").Value = "  -4.45%  "

$ws.Range("D49").Value = "2.086"
$ws.Range("E49").Value = "  -7.01%  "

$ws.Range("D50").Value = "1.265"
$ws.Range("E50").Value = "  +2.80%  "

$ws.Range("D51").Value = "124.15"
$ws.Range("E51").Value = "  -2.34%  "
